$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.696.69'
$ws.Range('E2').Value = '  +1.20%  '
$ws.Range('D3').Value = '2.613.95'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '600.89'
$ws.Range('E5').Value = '  +1.12%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.42'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.551'
$ws.Range('D9').Value = '2.612.77'
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('E10').Value = '  +10.51%  '
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.23'
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('E13').Value = '  -1.05%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.58'
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('E15').Value = '  +3.23%  '
$ws.Range('D16').Value = '3.090.48'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '67.662.09'
$ws.Range('E17').Value = '  +1.55%  '
$ws.Range('D18').Value = '2.616.37'
$ws.Range('E18').Value = '  +0.97%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.16'
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '365.20'
$ws.Range('E20').Value = '  +2.90%  '
$ws.Range('E21').Value = '  -1.78%  '
$ws.Range('E22').Value = '  -0.58%  '
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '70.25'
$ws.Range('E25').Value = '  +4.38%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.83'
$ws.Range('E26').Value = '  -6.63%  '
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '575.05'
$ws.Range('E29').Value = '  -4.23%  '
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('E31').Value = '  -2.50%  '
$ws.Range('E32').Value = '  -2.38%  '
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('E34').Value = '  -1.68%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('E36').Value = '  -3.65%  '
$ws.Range('E37').Value = '  -1.84%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '158.25'
$ws.Range('E38').Value = '  +2.96%  '
$ws.Range('E39').Value = '  +0.61%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('E41').Value = '  -2.38%  '
$ws.Range('E42').Value = '  +2.68%  '
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '41.16'
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('E45').Value = '  +0.09%  '
$ws.Range('E46').Value = '  -0.20%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '156.50'
$ws.Range('E47').Value = '  +0.50%  '
$ws.Range('D48').Value = '0.0₆0286'
$ws.Range('E48').Value = '  -7.60%  '
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('E50').Value = '  -2.16%  '
$ws.Range('E51').Value = '  +1.41%  '
